$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells being updated, to prevent Excel
# auto-converting numeric-looking strings (e.g. "1.00") into numbers.
$priceCells = @(2,3,5,6,7,9,10,12,14,15,16,17,18,20,21,22,23,25,26,27,29,30,31,32,34,35,36,37,38,40,41,42,44,45,46,49,50,51)
foreach ($r in $priceCells) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '67.335.83'
$ws.Range("E2").Value = '  -1.76%  '

# Row 3
$ws.Range("D3").Value = '3.751.30'
$ws.Range("E3").Value = '  -0.73%  '

# Row 4
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").Value = '592.03'
$ws.Range("E5").Value = '  -0.72%  '

# Row 6
$ws.Range("D6").Value = '165.03'
$ws.Range("E6").Value = '  -2.11%  '

# Row 7
$ws.Range("D7").Value = '3.749.89'
$ws.Range("E7").Value = '  -0.73%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  -0.91%  '

# Row 10
$ws.Range("D10").Value = '0.157'
$ws.Range("E10").Value = '  -3.14%  '

# Row 11
$ws.Range("E11").Value = '  -2.25%  '

# Row 12
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  -0.47%  '

# Row 13
$ws.Range("E13").Value = '  -4.11%  '

# Row 14
$ws.Range("D14").Value = '35.61'
$ws.Range("E14").Value = '  -2.89%  '

# Row 15
$ws.Range("D15").Value = '4.387.49'
$ws.Range("E15").Value = '  -0.55%  '

# Row 16
$ws.Range("D16").Value = '3.758.79'
$ws.Range("E16").Value = '  -0.43%  '

# Row 17
$ws.Range("D17").Value = '67.389.92'
$ws.Range("E17").Value = '  -1.59%  '

# Row 18
$ws.Range("D18").Value = '17.57'
$ws.Range("E18").Value = '  -4.00%  '

# Row 19
$ws.Range("E19").Value = '  +0.13%  '

# Row 20
$ws.Range("D20").Value = '6.88'
$ws.Range("E20").Value = '  -2.64%  '

# Row 21
$ws.Range("D21").Value = '10.45'
$ws.Range("E21").Value = '  -4.61%  '

# Row 22
$ws.Range("D22").Value = '453.98'
$ws.Range("E22").Value = '  -3.00%  '

# Row 23
$ws.Range("D23").Value = '0.691'
$ws.Range("E23").Value = '  -1.89%  '

# Row 24
$ws.Range("E24").Value = '  +4.49%  '

# Row 25
$ws.Range("D25").Value = '82.96'
$ws.Range("E25").Value = '  -2.43%  '

# Row 26
$ws.Range("D26").Value = '2.12'
$ws.Range("E26").Value = '  -5.30%  '

# Row 27
$ws.Range("D27").Value = '11.78'
$ws.Range("E27").Value = '  -3.51%  '

# Row 28
$ws.Range("E28").Value = '  -0.08%  '

# Row 29
$ws.Range("D29").Value = '9.98'
$ws.Range("E29").Value = '  -2.08%  '

# Row 30
$ws.Range("D30").Value = '2.76'
$ws.Range("E30").Value = '  -1.34%  '

# Row 31
$ws.Range("D31").Value = '29.56'
$ws.Range("E31").Value = '  -1.89%  '

# Row 32
$ws.Range("D32").Value = '7.15'
$ws.Range("E32").Value = '  -3.55%  '

# Row 33
$ws.Range("E33").Value = '  -3.56%  '

# Row 34
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.02%  '

# Row 35
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '9.09'
$ws.Range("E35").Value = '  -2.18%  '

# Row 36
$ws.Range("D36").Value = '3.706.75'
$ws.Range("E36").Value = '  -0.69%  '

# Row 37
$ws.Range("D37").Value = '0.0993'
$ws.Range("E37").Value = '  -2.40%  '

# Row 38
$ws.Range("D38").Value = '3.30'
$ws.Range("E38").Value = '  -4.48%  '

# Row 39
$ws.Range("E39").Value = '  -2.12%  '

# Row 40
$ws.Range("D40").Value = '0.989'
$ws.Range("E40").Value = '  -1.36%  '

# Row 41
$ws.Range("D41").Value = '5.70'
$ws.Range("E41").Value = '  -2.42%  '

# Row 42
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.16%  '

# Row 44
$ws.Range("D44").Value = '43.76'
$ws.Range("E44").Value = '  +0.34%  '

# Row 45
$ws.Range("D45").Value = '0.297'
$ws.Range("E45").Value = '  -3.65%  '

# Row 46
$ws.Range("D46").Value = '46.79'
$ws.Range("E46").Value = '  +2.33%  '

# Row 47
$ws.Range("E47").Value = '  -4.73%  '

# Row 48
$ws.Range("E48").Value = '  -3.22%  '

# Row 49
$ws.Range("D49").Value = '146.51'
$ws.Range("E49").Value = '  +0.66%  '

# Row 50
$ws.Range("D50").Value = '388.74'
$ws.Range("E50").Value = '  -4.90%  '

# Row 51
$ws.Range("D51").Value = '2.740.35'
$ws.Range("E51").Value = '  +1.96%  '

# Restore default style on the Price cells we touched (clears the temporary
# text number-format so styles.xml matches the original formatting).
foreach ($r in $priceCells) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}